$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOCO")

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy cell formatting from the shifted columns (F:G, which now hold what used
# to be in D:E) into the freshly inserted D:E columns so the new quarter columns
# pick up the same number formats / styles as their neighbours. Restricted to
# the three statement blocks that actually carry data (skipping the blank
# spacer/header rows) so we don't manufacture empty styled cells where none
# existed before.
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarterly columns (D = most recent quarter, E = previous
# quarter) with the newly reported figures.
$newQuarterData = @{
    7 = @(43460, 43369)
    8 = @(106300, 112200)
    9 = @(54600, 57500)
    10 = @(51700, 54700)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(37000, -2200)
    15 = @(4800, 4500)
    17 = @(137300, 102700)
    18 = @(-31000, 9500)
    20 = @(-900, -300)
    21 = @(-27100, 13700)
    22 = @(0, 0)
    23 = @(-31800, 9200)
    24 = @(-8400, 2400)
    25 = @(0, 0)
    26 = @(-23400, 6800)
    27 = @(-23400, 6800)
    28 = @(0, 0)
    29 = @(0, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(900, 300)
    33 = @(-23400, 6800)
    34 = @(0, 0)
    35 = @(-23400, 6800)
    38 = @(43460, 43369)
    41 = @(7000, 8100)
    42 = @(0, 0)
    43 = @(9600, 8900)
    44 = @(2500, 2200)
    45 = @(3000, 2500)
    46 = @(22000, 21700)
    47 = @(0, 0)
    48 = @(104200, 105200)
    49 = @(310800, 310900)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(13200, 7100)
    53 = @(0, 0)
    54 = @(450200, 444900)
    57 = @(9600, 6800)
    58 = @(100, 100)
    59 = @(73300, 44900)
    60 = @(82900, 51900)
    61 = @(74100, 71100)
    62 = @(28000, 32900)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(185000, 155900)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-110900, -87500)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(265200, 289000)
    77 = @(0, 0)
    80 = @(43460, 43369)
    81 = @(-23400, 6800)
    83 = @(4800, 4500)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(3800, 16000)
    91 = @(-7100, -6800)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-7100, -6800)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(2100, -14200)
    101 = @(0, 0)
    102 = @(-1200, -4900)
}

foreach ($r in $newQuarterData.Keys) {
    $vals = $newQuarterData[$r]
    $ws.Cells.Item($r, 4).Value2 = $vals[0]
    $ws.Cells.Item($r, 5).Value2 = $vals[1]
}

Write-Host "Inserted 2 columns and populated new quarter data for $($newQuarterData.Count) rows."
